$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.277.44'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.38%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.832.32'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.17%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.88%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9992'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.12%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4977'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -3.09%  '

$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3936'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.58%  '

$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1004'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +27.68%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.112'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.29%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.11'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.20%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.454'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.15%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.70'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.45%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.9997'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.09%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.823.41'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.85%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.343'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.17%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001144'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +5.61%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '93.04'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.25%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06645'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.04%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9988'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.13%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.26'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.35%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.042'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.40%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.308.87'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.45%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.50%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.233'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.00%  '

$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.81'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.34%  '

$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '158.10'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.32%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.039.37'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.94%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.437'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.01%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.16'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.85%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1054'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.02%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.047'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.06%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.609'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.51%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.605'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.63%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.06785'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -6.48%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.086'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.78%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02352'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.35%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2154'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.97%  '

$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '11.47'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.81%  '

$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.993'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.28%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6235'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.55%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.180'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.95%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9989'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.14%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.26'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.41%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5949'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.72%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.687'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.36%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.271'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.37%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.22%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.952'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.02%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.184'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.14%  '

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'ThetaToken'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.118'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +4.22%  '
